# Scheduled price/profit refresh for the Leve profit sheets.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H, I, J, K, L, M, N)
# for the rows whose market data changed, across the ALC, ARM, BSM, CRP, CUL,
# GSM, LTW and WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 429347.03
$ws.Range("J17").Value = 490553.75
$ws.Range("L17").Value = 1471661.25
$ws.Range("N17").Value = -1471997.25
$ws.Range("H40").Value = 3335880.5
$ws.Range("I40").Value = 2939.8
$ws.Range("K40").Value = 2939.8
$ws.Range("M40").Value = -2764.8
$ws.Range("H132").Value = 1643
$ws.Range("J132").Value = 2003.6666
$ws.Range("L132").Value = 6010.9998
$ws.Range("N132").Value = -11070.9998
$ws.Range("H138").Value = 4720.9736
$ws.Range("I138").Value = 1658.091
$ws.Range("J138").Value = 8932.4375
$ws.Range("K138").Value = 4974.272999999999
$ws.Range("L138").Value = 26797.3125
$ws.Range("M138").Value = 165.7270000000008
$ws.Range("N138").Value = -37077.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5442.222
$ws.Range("I45").Value = 1620.9231
$ws.Range("K45").Value = 1620.9231
$ws.Range("M45").Value = -1243.9231
$ws.Range("H132").Value = 2945.8333
$ws.Range("I132").Value = 1744.2258
$ws.Range("J132").Value = 6332.1816
$ws.Range("K132").Value = 5232.6774
$ws.Range("L132").Value = 18996.5448
$ws.Range("M132").Value = -2702.6774
$ws.Range("N132").Value = -24056.5448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8773803
$ws.Range("J20").Value = 1693.7273
$ws.Range("L20").Value = 1693.7273
$ws.Range("N20").Value = -2187.7273
$ws.Range("H86").Value = 37078520
$ws.Range("I86").Value = 62497.35
$ws.Range("K86").Value = 62497.35
$ws.Range("M86").Value = -61374.35
$ws.Range("H89").Value = 37078520
$ws.Range("I89").Value = 62497.35
$ws.Range("K89").Value = 312486.75
$ws.Range("M89").Value = -306870.75
$ws.Range("H107").Value = 32146068
$ws.Range("I107").Value = 48915332
$ws.Range("K107").Value = 48915332
$ws.Range("M107").Value = -48913412
$ws.Range("H134").Value = 4634.648
$ws.Range("I134").Value = 1217.5
$ws.Range("K134").Value = 3652.5
$ws.Range("M134").Value = -1117.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7151.413
$ws.Range("I31").Value = 2543.7083
$ws.Range("K31").Value = 2543.7083
$ws.Range("M31").Value = -2248.7083
$ws.Range("H34").Value = 7151.413
$ws.Range("I34").Value = 2543.7083
$ws.Range("K34").Value = 2543.7083
$ws.Range("M34").Value = -2341.7083
$ws.Range("H86").Value = 7880798.5
$ws.Range("J86").Value = 129500
$ws.Range("L86").Value = 129500
$ws.Range("N86").Value = -131746
$ws.Range("H89").Value = 7880798.5
$ws.Range("J89").Value = 129500
$ws.Range("L89").Value = 647500
$ws.Range("N89").Value = -658732
$ws.Range("H99").Value = 10262.308
$ws.Range("I99").Value = 12067.667
$ws.Range("K99").Value = 12067.667
$ws.Range("M99").Value = -10569.667
$ws.Range("H126").Value = 10262.308
$ws.Range("I126").Value = 12067.667
$ws.Range("K126").Value = 36203.001
$ws.Range("M126").Value = -33733.001
$ws.Range("H132").Value = 4563.553
$ws.Range("I132").Value = 1936
$ws.Range("K132").Value = 5808
$ws.Range("M132").Value = -3278
$ws.Range("H134").Value = 7187.7334
$ws.Range("I134").Value = 1768
$ws.Range("K134").Value = 5304
$ws.Range("M134").Value = -2769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 202462.8
$ws.Range("I140").Value = 287232.56
$ws.Range("K140").Value = 861697.6799999999
$ws.Range("M140").Value = -856517.6799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 69999.664
$ws.Range("J52").Value = 89999.5
$ws.Range("L52").Value = 89999.5
$ws.Range("N52").Value = -90517.5
$ws.Range("H80").Value = 3490
$ws.Range("I80").Value = 3735.5
$ws.Range("J80").Value = 2999
$ws.Range("K80").Value = 3735.5
$ws.Range("L80").Value = 2999
$ws.Range("M80").Value = -2737.5
$ws.Range("N80").Value = -4995
$ws.Range("H83").Value = 3490
$ws.Range("I83").Value = 3735.5
$ws.Range("J83").Value = 2999
$ws.Range("K83").Value = 18677.5
$ws.Range("L83").Value = 14995
$ws.Range("M83").Value = -13685.5
$ws.Range("N83").Value = -24979
$ws.Range("H107").Value = 800440
$ws.Range("I107").Value = 1143128.6
$ws.Range("K107").Value = 1143128.6
$ws.Range("M107").Value = -1141208.6
$ws.Range("H113").Value = 7079.3105
$ws.Range("J113").Value = 9322.223
$ws.Range("L113").Value = 9322.223
$ws.Range("N113").Value = -13662.223
$ws.Range("H132").Value = 1682.3226
$ws.Range("I132").Value = 1511.3928
$ws.Range("J132").Value = 3277.6667
$ws.Range("K132").Value = 4534.178400000001
$ws.Range("L132").Value = 9833.000100000001
$ws.Range("M132").Value = -2004.178400000001
$ws.Range("N132").Value = -14893.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 860
$ws.Range("I16").Value = 575
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 575
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -405
$ws.Range("N16").Value = -2340
$ws.Range("H22").Value = 1700.0555
$ws.Range("I22").Value = 919.2308
$ws.Range("K22").Value = 919.2308
$ws.Range("M22").Value = -624.2308
$ws.Range("H27").Value = 1700.0555
$ws.Range("I27").Value = 919.2308
$ws.Range("K27").Value = 919.2308
$ws.Range("M27").Value = -812.2308
$ws.Range("H29").Value = 1137
$ws.Range("I29").Value = 1137
$ws.Range("K29").Value = 1137
$ws.Range("M29").Value = -842
$ws.Range("H39").Value = 2000
$ws.Range("I39").Value = 2000
$ws.Range("K39").Value = 2000
$ws.Range("M39").Value = -1540
$ws.Range("H46").Value = 2301534.2
$ws.Range("I46").Value = 34482760
$ws.Range("J46").Value = 2875.3572
$ws.Range("K46").Value = 34482760
$ws.Range("L46").Value = 2875.3572
$ws.Range("M46").Value = -34482572
$ws.Range("N46").Value = -3251.3572
$ws.Range("H61").Value = 4872.5713
$ws.Range("I61").Value = 907.4
$ws.Range("K61").Value = 907.4
$ws.Range("M61").Value = -705.4
$ws.Range("H93").Value = 5648.1333
$ws.Range("I93").Value = 3657.2222
$ws.Range("J93").Value = 8634.5
$ws.Range("K93").Value = 3657.2222
$ws.Range("L93").Value = 8634.5
$ws.Range("M93").Value = -2409.2222
$ws.Range("N93").Value = -11130.5
$ws.Range("H113").Value = 4872.5713
$ws.Range("I113").Value = 907.4
$ws.Range("K113").Value = 907.4
$ws.Range("M113").Value = 1262.6
$ws.Range("H122").Value = 3769.6
$ws.Range("I122").Value = 2793.6428
$ws.Range("K122").Value = 8380.928400000001
$ws.Range("M122").Value = -5930.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14718127
$ws.Range("I132").Value = 18523312
$ws.Range("J132").Value = 40983.855
$ws.Range("K132").Value = 55569936
$ws.Range("L132").Value = 122951.565
$ws.Range("M132").Value = -55567406
$ws.Range("N132").Value = -128011.565
